$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing rows 52-79 (date shifted one market-day block earlier; volumen/precio/unidad/origen refreshed)
$ws.Range("D52").Value = 44574
$ws.Range("I52").Value = 'Extra'
$ws.Range("J52").Value = 400
$ws.Range("K52").Value = 3000
$ws.Range("L52").Value = 3000
$ws.Range("M52").Value = 3000
$ws.Range("N52").Value = '$/unidad'
$ws.Range("O52").Value = 'Región de O''Higgins'
$ws.Range("P52").Value = 3000

$ws.Range("D53").Value = 44574
$ws.Range("I53").Value = 'Primera'
$ws.Range("J53").Value = 400
$ws.Range("K53").Value = 2500
$ws.Range("L53").Value = 2500
$ws.Range("M53").Value = 2500
$ws.Range("N53").Value = '$/unidad'
$ws.Range("O53").Value = 'Región de O''Higgins'
$ws.Range("P53").Value = 2500

$ws.Range("D54").Value = 44574
$ws.Range("I54").Value = 'Segunda'
$ws.Range("J54").Value = 400
$ws.Range("K54").Value = 2000
$ws.Range("L54").Value = 2000
$ws.Range("M54").Value = 2000
$ws.Range("N54").Value = '$/unidad'
$ws.Range("O54").Value = 'Región de O''Higgins'
$ws.Range("P54").Value = 2000

$ws.Range("D55").Value = 44251
$ws.Range("I55").Value = 'Extra'
$ws.Range("J55").Value = 300
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 3000
$ws.Range("M55").Value = 3000
$ws.Range("N55").Value = '$/unidad'
$ws.Range("O55").Value = 'Región de O''Higgins'
$ws.Range("P55").Value = 3000

$ws.Range("D56").Value = 44251
$ws.Range("I56").Value = 'Primera'
$ws.Range("J56").Value = 500
$ws.Range("K56").Value = 2500
$ws.Range("L56").Value = 2500
$ws.Range("M56").Value = 2500
$ws.Range("N56").Value = '$/unidad'
$ws.Range("O56").Value = 'Región de O''Higgins'
$ws.Range("P56").Value = 2500

$ws.Range("D57").Value = 44251
$ws.Range("I57").Value = 'Segunda'
$ws.Range("J57").Value = 500
$ws.Range("K57").Value = 2000
$ws.Range("L57").Value = 2000
$ws.Range("M57").Value = 2000
$ws.Range("N57").Value = '$/unidad'
$ws.Range("O57").Value = 'Región de O''Higgins'
$ws.Range("P57").Value = 2000

$ws.Range("D58").Value = 44272
$ws.Range("I58").Value = 'Primera'
$ws.Range("J58").Value = 300
$ws.Range("K58").Value = 2500
$ws.Range("L58").Value = 2500
$ws.Range("M58").Value = 2500
$ws.Range("N58").Value = '$/unidad'
$ws.Range("O58").Value = 'Región de O''Higgins'
$ws.Range("P58").Value = 2500

$ws.Range("D59").Value = 44272
$ws.Range("I59").Value = 'Segunda'
$ws.Range("J59").Value = 300
$ws.Range("K59").Value = 2000
$ws.Range("L59").Value = 2000
$ws.Range("M59").Value = 2000
$ws.Range("N59").Value = '$/unidad'
$ws.Range("O59").Value = 'Región de O''Higgins'
$ws.Range("P59").Value = 2000

$ws.Range("D60").Value = 44211
$ws.Range("I60").Value = 'Extra'
$ws.Range("J60").Value = 500
$ws.Range("K60").Value = 3500
$ws.Range("L60").Value = 3500
$ws.Range("M60").Value = 3500
$ws.Range("N60").Value = '$/unidad'
$ws.Range("O60").Value = 'Región de O''Higgins'
$ws.Range("P60").Value = 3500

$ws.Range("D61").Value = 44211
$ws.Range("I61").Value = 'Primera'
$ws.Range("J61").Value = 500
$ws.Range("K61").Value = 3000
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = 3000
$ws.Range("N61").Value = '$/unidad'
$ws.Range("O61").Value = 'Región de O''Higgins'
$ws.Range("P61").Value = 3000

$ws.Range("D62").Value = 44211
$ws.Range("I62").Value = 'Segunda'
$ws.Range("J62").Value = 500
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 2500
$ws.Range("M62").Value = 2500
$ws.Range("N62").Value = '$/unidad'
$ws.Range("O62").Value = 'Región de O''Higgins'
$ws.Range("P62").Value = 2500

$ws.Range("D63").Value = 44196
$ws.Range("I63").Value = 'Extra'
$ws.Range("J63").Value = 400
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = 3000
$ws.Range("N63").Value = '$/unidad'
$ws.Range("O63").Value = 'Región de O''Higgins'
$ws.Range("P63").Value = 3000

$ws.Range("D64").Value = 44196
$ws.Range("I64").Value = 'Primera'
$ws.Range("J64").Value = 400
$ws.Range("K64").Value = 2500
$ws.Range("L64").Value = 2500
$ws.Range("M64").Value = 2500
$ws.Range("N64").Value = '$/unidad'
$ws.Range("O64").Value = 'Región de O''Higgins'
$ws.Range("P64").Value = 2500

$ws.Range("D65").Value = 44196
$ws.Range("I65").Value = 'Segunda'
$ws.Range("J65").Value = 400
$ws.Range("K65").Value = 2000
$ws.Range("L65").Value = 2000
$ws.Range("M65").Value = 2000
$ws.Range("N65").Value = '$/unidad'
$ws.Range("O65").Value = 'Región de O''Higgins'
$ws.Range("P65").Value = 2000

$ws.Range("D66").Value = 44208
$ws.Range("I66").Value = 'Extra'
$ws.Range("J66").Value = 500
$ws.Range("K66").Value = 3500
$ws.Range("L66").Value = 3500
$ws.Range("M66").Value = 3500
$ws.Range("N66").Value = '$/kilo (volumen en unidades)'
$ws.Range("O66").Value = 'Región de O''Higgins'
$ws.Range("P66").Value = 3500

$ws.Range("D67").Value = 44208
$ws.Range("I67").Value = 'Primera'
$ws.Range("J67").Value = 500
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 3000
$ws.Range("M67").Value = 3000
$ws.Range("N67").Value = '$/kilo (volumen en unidades)'
$ws.Range("O67").Value = 'Región de O''Higgins'
$ws.Range("P67").Value = 3000

$ws.Range("D68").Value = 44208
$ws.Range("I68").Value = 'Segunda'
$ws.Range("J68").Value = 500
$ws.Range("K68").Value = 2500
$ws.Range("L68").Value = 2500
$ws.Range("M68").Value = 2500
$ws.Range("N68").Value = '$/kilo (volumen en unidades)'
$ws.Range("O68").Value = 'Región de O''Higgins'
$ws.Range("P68").Value = 2500

$ws.Range("D69").Value = 44545
$ws.Range("I69").Value = 'Primera'
$ws.Range("J69").Value = 1100
$ws.Range("K69").Value = 2500
$ws.Range("L69").Value = 3000
$ws.Range("M69").Value = 2773
$ws.Range("N69").Value = '$/unidad'
$ws.Range("O69").Value = 'Región de O''Higgins'
$ws.Range("P69").Value = 2773

$ws.Range("D70").Value = 44545
$ws.Range("I70").Value = 'Segunda'
$ws.Range("J70").Value = 350
$ws.Range("K70").Value = 2300
$ws.Range("L70").Value = 2500
$ws.Range("M70").Value = 2414
$ws.Range("N70").Value = '$/unidad'
$ws.Range("O70").Value = 'Región de O''Higgins'
$ws.Range("P70").Value = 2414

$ws.Range("D71").Value = 44194
$ws.Range("I71").Value = 'Extra'
$ws.Range("J71").Value = 400
$ws.Range("K71").Value = 3000
$ws.Range("L71").Value = 3000
$ws.Range("M71").Value = 3000
$ws.Range("N71").Value = '$/unidad'
$ws.Range("O71").Value = 'Región de O''Higgins'
$ws.Range("P71").Value = 3000

$ws.Range("D72").Value = 44194
$ws.Range("I72").Value = 'Primera'
$ws.Range("J72").Value = 500
$ws.Range("K72").Value = 2500
$ws.Range("L72").Value = 2500
$ws.Range("M72").Value = 2500
$ws.Range("N72").Value = '$/unidad'
$ws.Range("O72").Value = 'Región de O''Higgins'
$ws.Range("P72").Value = 2500

$ws.Range("D73").Value = 44194
$ws.Range("I73").Value = 'Segunda'
$ws.Range("J73").Value = 400
$ws.Range("K73").Value = 2000
$ws.Range("L73").Value = 2000
$ws.Range("M73").Value = 2000
$ws.Range("N73").Value = '$/unidad'
$ws.Range("O73").Value = 'Región de O''Higgins'
$ws.Range("P73").Value = 2000

$ws.Range("D74").Value = 44236
$ws.Range("I74").Value = 'Extra'
$ws.Range("J74").Value = 500
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = 3000
$ws.Range("N74").Value = '$/unidad'
$ws.Range("O74").Value = 'Región de O''Higgins'
$ws.Range("P74").Value = 3000

$ws.Range("D75").Value = 44236
$ws.Range("I75").Value = 'Primera'
$ws.Range("J75").Value = 500
$ws.Range("K75").Value = 2500
$ws.Range("L75").Value = 2500
$ws.Range("M75").Value = 2500
$ws.Range("N75").Value = '$/unidad'
$ws.Range("O75").Value = 'Región de O''Higgins'
$ws.Range("P75").Value = 2500

$ws.Range("D76").Value = 44236
$ws.Range("I76").Value = 'Segunda'
$ws.Range("J76").Value = 500
$ws.Range("K76").Value = 2000
$ws.Range("L76").Value = 2000
$ws.Range("M76").Value = 2000
$ws.Range("N76").Value = '$/unidad'
$ws.Range("O76").Value = 'Región de O''Higgins'
$ws.Range("P76").Value = 2000

$ws.Range("D77").Value = 44264
$ws.Range("I77").Value = 'Extra'
$ws.Range("J77").Value = 300
$ws.Range("K77").Value = 2800
$ws.Range("L77").Value = 2800
$ws.Range("M77").Value = 2800
$ws.Range("N77").Value = '$/unidad'
$ws.Range("O77").Value = 'Región de O''Higgins'
$ws.Range("P77").Value = 2800

$ws.Range("D78").Value = 44264
$ws.Range("I78").Value = 'Primera'
$ws.Range("J78").Value = 300
$ws.Range("K78").Value = 2500
$ws.Range("L78").Value = 2500
$ws.Range("M78").Value = 2500
$ws.Range("N78").Value = '$/unidad'
$ws.Range("O78").Value = 'Región de O''Higgins'
$ws.Range("P78").Value = 2500

$ws.Range("D79").Value = 44264
$ws.Range("I79").Value = 'Segunda'
$ws.Range("J79").Value = 300
$ws.Range("K79").Value = 2200
$ws.Range("L79").Value = 2200
$ws.Range("M79").Value = 2200
$ws.Range("N79").Value = '$/unidad'
$ws.Range("O79").Value = 'Región de O''Higgins'
$ws.Range("P79").Value = 2200

# Append two new rows (80, 81) with data that used to belong to Peru-origin entries
$ws.Range("A80").Value = 11
$ws.Range("B80").Value = 'Vega Monumental Concepción'
$ws.Range("C80").Value = 'Bíobío'
$ws.Range("D80").Value = 44525
$ws.Range("E80").Value = 8
$ws.Range("F80").Value = 100112028
$ws.Range("G80").Value = 'Sandia'
$ws.Range("H80").Value = 'Sin especificar'
$ws.Range("I80").Value = 'Primera'
$ws.Range("J80").Value = 200
$ws.Range("K80").Value = 700
$ws.Range("L80").Value = 800
$ws.Range("M80").Value = 750
$ws.Range("N80").Value = '$/kilo (volumen en unidades)'
$ws.Range("O80").Value = 'Perú'
$ws.Range("P80").Value = 750
$ws.Range("Q80").Value = 1
$ws.Range("R80").Value = 'Hortaliza'

$ws.Range("A81").Value = 11
$ws.Range("B81").Value = 'Vega Monumental Concepción'
$ws.Range("C81").Value = 'Bíobío'
$ws.Range("D81").Value = 44511
$ws.Range("E81").Value = 8
$ws.Range("F81").Value = 100112028
$ws.Range("G81").Value = 'Sandia'
$ws.Range("H81").Value = 'Sin especificar'
$ws.Range("I81").Value = 'Primera'
$ws.Range("J81").Value = 600
$ws.Range("K81").Value = 800
$ws.Range("L81").Value = 900
$ws.Range("M81").Value = 850
$ws.Range("N81").Value = '$/kilo (volumen en unidades)'
$ws.Range("O81").Value = 'Perú'
$ws.Range("P81").Value = 850
$ws.Range("Q81").Value = 1
$ws.Range("R81").Value = 'Hortaliza'
